$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "cid"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "name_kor"

# Data rows: cid, name (code), name_kor, description, group
$rows = @(
  @(1, "WalkersWorkshop",  "워커 공방",           "hp 10% 증가",                      "weak"),
  @(2, "SmithsSmithy",     "스미스 제철소",        "attack damage 10% 증가",           "weak"),
  @(3, "DSDC",             "슬럼가 개발 주식회사",  "예측할 수 없는 특수효과",           "weak"),
  @(4, "KanaKooler",       "카나 쿨러",            "쿨타임 10% 감소",                  "gang"),
  @(5, "RaccoonGlass",     "라쿤 유리 공예",        "피격시 10% 확률로 데미지 반사",     "gang"),
  @(6, "DelicateBastard",  "섬세한 불한당",         "피흡 10%",                         "fight"),
  @(7, "ProIndustrial",    "프로 공업소",           "10% 확률로 적 마비",                "fight"),
  @(8, "RoyalLab",         "왕국 연구소",           "스킬 10회 사용시 Divine 1회 적용",  "master")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 7).Value = $row[4]
    $r = $r + 1
}

# Column D width
$ws.Columns.Item(4).ColumnWidth = 15.41

# Selection
$ws.Range("D9").Select() | Out-Null

# Sheet-scoped hidden defined name for the filter database (mirrors AutoFilter metadata)
$fd = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$G`$1")
$fd.Visible = $false
